$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# The update adds a brand-new match result row ("Odisha FC" vs
# "Kerala Blasters", id 132) that sorts by date BEFORE the previous
# last row ("FC Goa" vs "Chennaiyin FC", id 132 -> becomes id 133).
# So: push the old last row (134) down to row 135 (re-numbering the
# id and refreshing two odds values), then overwrite row 134 with the
# freshly finished match.
# -----------------------------------------------------------------

# --- Step 1: move the previous last row (134) down to row 135 -----
# Copy formats first (from cells that already carry the exact same
# style index) so no new cellXfs entries get minted.
$ws.Range("A133").Copy() | Out-Null
$ws.Range("A135").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E133").Copy() | Out-Null
$ws.Range("E135").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A135").Value = 133
$ws.Range("B135").Value = 8103574
$ws.Range("C135").Value = "India Super League"
$ws.Range("D135").Value = "India Super League"
$ws.Range("E135").Value = 45402.45833333334
$ws.Range("F135").Value = "FC Goa"
$ws.Range("G135").Value = "Chennaiyin FC"
$ws.Range("K135").Value = 1.4
$ws.Range("L135").Value = 4.5
$ws.Range("M135").Value = 6
$ws.Range("N135").Value = 1.4
$ws.Range("O135").Value = 4.5
$ws.Range("P135").Value = 5.75
$ws.Range("Q135").Value = -1.25
$ws.Range("R135").Value = 1.925
$ws.Range("S135").Value = 1.875
$ws.Range("T135").Value = 3
$ws.Range("U135").Value = 1.825
$ws.Range("V135").Value = 1.975
$ws.Range("W135").Value = 0
$ws.Range("X135").Value = 0
$ws.Range("Y135").Value = 0
$ws.Range("Z135").Value = 0
$ws.Range("AA135").Value = 0

# --- Step 2: overwrite row 134 with the new finished match --------
$ws.Range("A134").Value = 132
$ws.Range("B134").Value = 8103573
$ws.Range("C134").Value = "India Super League"
$ws.Range("D134").Value = "India Super League"
$ws.Range("E134").Value = 45401.45833333334
$ws.Range("F134").Value = "Odisha FC"
$ws.Range("G134").Value = "Kerala Blasters"
$ws.Range("H134").Value = 2
$ws.Range("I134").Value = 1
$ws.Range("J134").Value = "H"
$ws.Range("K134").Value = 1.65
$ws.Range("L134").Value = 3.7
$ws.Range("M134").Value = 4.5
$ws.Range("N134").Value = 1.75
$ws.Range("O134").Value = 3.5
$ws.Range("P134").Value = 4.2
$ws.Range("Q134").Value = -0.75
$ws.Range("R134").Value = 2.025
$ws.Range("S134").Value = 1.825
$ws.Range("T134").Value = 2.75
$ws.Range("U134").Value = 2
$ws.Range("V134").Value = 1.85
$ws.Range("W134").Value = 0.75
$ws.Range("X134").Value = -1
$ws.Range("Y134").Value = -1
$ws.Range("Z134").Value = 0.5125
$ws.Range("AA134").Value = -0.5
$ws.Range("AB134").Value = 0.5
$ws.Range("AC134").Value = -0.5
